# Apply edits described by the commit:
# "Fixed test_dynamic_components.py :: TestInitialRate_dM test case; incorporated
#  1/3600 h/s factor into test_model_for_exponential_growth_in_M.xlsx model
#  growthRate parameter."

$wb = $excel.ActiveWorkbook

$wsRateLaws   = $wb.Worksheets.Item("Rate laws")
$wsParameters = $wb.Worksheets.Item("Parameters")

# --- "Rate laws" sheet: C2 formula text loses the /3600 factor, since that
#     factor now lives inside the growthRate parameter itself.
$wsRateLaws.Range("C2").Value = "growthRate * M[c]"

# --- "Parameters" sheet: growthRate row (row 3)
#     E3 numeric value updated to incorporate the 1/3600 h/s factor, and
#     formatted with a scientific number format (existing style).
$wsParameters.Range("E3").Value = 0.0000083713
$wsParameters.Range("E3").NumberFormat = "0.00E+00"

#     G3 comment explains how the new value was derived.
$wsParameters.Range("G3").Value = "ln(2)/23 h * 1 h / 3600 s = ln(2)/(23*3600) 1/s = 8.3713e-06 1/s"

# --- Selections / active sheet bookkeeping: author had moved on from the
#     "Rate laws" tab to the "Parameters" tab.
$wsRateLaws.Range("A3").Select()
$wsParameters.Range("A4").Select()

$wsParameters.Activate()
